# feat(house): implement upgrade panel
# Adds the localization rows (key / en / fr) for the new upgrade panel UI
# to the end of the translations table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("UPGRADE_TITLE",       "UPGRADE", "AMÉLIORER"),
    @("UPGRADE_STONE_LABEL", "Stone",   "Pierre"),
    @("UPGRADE_GOLD_LABEL",  "Gold",    "Or"),
    @("UPGRADE_BUTTON",      "UPGRADE", "AMÉLIORER")
)

$r = 19
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("F20").Select()
